$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")

# Insert a new row at row 2 (pushes existing data rows down by one)
$ws.Rows.Item(2).Insert()

# Populate the new order row (force text format on text-like columns so
# numeric/date-looking strings such as the phone number and collection
# date are not auto-coerced to numbers/dates)
$ws.Range("A2").Value = 4
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2026-01-13 16:39"
$ws.Range("C2").Value = "Pooja"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "A1608"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "9096648553"
$ws.Range("F2").Value = "Onion Pakoda (Kanda Bhaje) x1"
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2026-01-14"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "22:09"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

# Update the Daily Summary sheet totals
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Range("B2").Value = 4
$ws2.Range("E2").Value = 135
$ws2.Range("G2").Value = 135
